$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18 (id 5471)
$ws.Range("H18").Value = 688.3333
$ws.Range("I18").Value = 688.3333
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 688.3333
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -404.3333
$ws.Range("N18").Value = ""

# Row 70 (id 12604)
$ws.Range("H70").Value = 2724.6287
$ws.Range("I70").Value = 5936.2
$ws.Range("J70").Value = 1440
$ws.Range("K70").Value = 17808.6
$ws.Range("L70").Value = 4320
$ws.Range("M70").Value = -17538.6
$ws.Range("N70").Value = -4860

# Row 73 (id 12604)
$ws.Range("H73").Value = 2724.6287
$ws.Range("I73").Value = 5936.2
$ws.Range("J73").Value = 1440
$ws.Range("K73").Value = 17808.6
$ws.Range("L73").Value = 4320
$ws.Range("M73").Value = -16872.6
$ws.Range("N73").Value = -6192

# Row 98 (id 36237)
$ws.Range("H98").Value = 828
$ws.Range("I98").Value = 596.36365
$ws.Range("J98").Value = 1394.2222
$ws.Range("K98").Value = 596.36365
$ws.Range("L98").Value = 1394.2222
$ws.Range("M98").Value = 901.63635
$ws.Range("N98").Value = -4390.2222

# Row 107 (id 27766)
$ws.Range("H107").Value = 373
$ws.Range("I107").Value = 267.55554
$ws.Range("J107").Value = 847.5
$ws.Range("K107").Value = 267.55554
$ws.Range("L107").Value = 847.5
$ws.Range("M107").Value = 1652.44446
$ws.Range("N107").Value = -4687.5

# Row 122 (id 36237)
$ws.Range("H122").Value = 828
$ws.Range("I122").Value = 596.36365
$ws.Range("J122").Value = 1394.2222
$ws.Range("K122").Value = 1789.09095
$ws.Range("L122").Value = 4182.6666
$ws.Range("M122").Value = 660.90905
$ws.Range("N122").Value = -9082.6666

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (id 44147)
$ws.Range("H32").Value = 7028.8096
$ws.Range("I32").Value = 4415.4326
$ws.Range("J32").Value = 26367.8
$ws.Range("K32").Value = 4415.4326
$ws.Range("L32").Value = 26367.8
$ws.Range("M32").Value = -4128.4326
$ws.Range("N32").Value = -26941.8

# Row 47 (id 3622)
$ws.Range("H47").Value = 16633.334
$ws.Range("J47").Value = 16633.334
$ws.Range("L47").Value = 16633.334
$ws.Range("N47").Value = -18083.334

# Row 110 (id 27708)
$ws.Range("H110").Value = 2390.818
$ws.Range("I110").Value = 2966.6667
$ws.Range("J110").Value = 1699.8
$ws.Range("K110").Value = 2966.6667
$ws.Range("L110").Value = 1699.8
$ws.Range("M110").Value = -921.6667000000002
$ws.Range("N110").Value = -5789.8

$ws = $wb.Worksheets.Item("BSM")
# Row 48 (id 22888)
$ws.Range("H48").Value = 69000
$ws.Range("J48").Value = 69000
$ws.Range("L48").Value = 69000
$ws.Range("N48").Value = -69830

# Row 105 (id 19947)
$ws.Range("H105").Value = 2123.3572
$ws.Range("I105").Value = 1974.0625
$ws.Range("J105").Value = 2601.1
$ws.Range("K105").Value = 1974.0625
$ws.Range("L105").Value = 2601.1
$ws.Range("M105").Value = -227.0625
$ws.Range("N105").Value = -6095.1

$ws = $wb.Worksheets.Item("CRP")
# Row 109 (id 27203)
$ws.Range("H109").Value = 28000
$ws.Range("J109").Value = 28000
$ws.Range("L109").Value = 28000
$ws.Range("N109").Value = -30080

$ws = $wb.Worksheets.Item("CUL")
# Row 131 (id 36060)
$ws.Range("H131").Value = 891.2
$ws.Range("I131").Value = 512.6667
$ws.Range("J131").Value = 909.5161000000001
$ws.Range("K131").Value = 1538.0001
$ws.Range("L131").Value = 2728.5483
$ws.Range("M131").Value = 3501.9999
$ws.Range("N131").Value = -12808.5483

$ws = $wb.Worksheets.Item("GSM")
# Row 108 (id 27082)
$ws.Range("H108").Value = 30000
$ws.Range("J108").Value = 30000
$ws.Range("L108").Value = 30000
$ws.Range("N108").Value = -37680

# Row 126 (id 36184)
$ws.Range("H126").Value = 2379.875
$ws.Range("I126").Value = 1627.8
$ws.Range("J126").Value = 3633.3333
$ws.Range("K126").Value = 4883.4
$ws.Range("L126").Value = 10899.9999
$ws.Range("M126").Value = -2413.4
$ws.Range("N126").Value = -15839.9999

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (id 36249)
$ws.Range("H7").Value = 2219.6538
$ws.Range("I7").Value = 2054.1333
$ws.Range("J7").Value = 2445.3635
$ws.Range("K7").Value = 2054.1333
$ws.Range("L7").Value = 2445.3635
$ws.Range("M7").Value = -1942.1333
$ws.Range("N7").Value = -2669.3635

# Row 22 (id 5277)
$ws.Range("H22").Value = 725.55554
$ws.Range("I22").Value = 735.7143
$ws.Range("J22").Value = 690
$ws.Range("K22").Value = 735.7143
$ws.Range("L22").Value = 690
$ws.Range("M22").Value = -440.7143
$ws.Range("N22").Value = -1280

# Row 27 (id 5277)
$ws.Range("H27").Value = 725.55554
$ws.Range("I27").Value = 735.7143
$ws.Range("J27").Value = 690
$ws.Range("K27").Value = 735.7143
$ws.Range("L27").Value = 690
$ws.Range("M27").Value = -628.7143
$ws.Range("N27").Value = -904

# Row 40 (id 36248)
$ws.Range("H40").Value = 2241.889
$ws.Range("I40").Value = 2208.9333
$ws.Range("J40").Value = 2406.6667
$ws.Range("K40").Value = 2208.9333
$ws.Range("L40").Value = 2406.6667
$ws.Range("M40").Value = -2072.9333
$ws.Range("N40").Value = -2678.6667

# Row 46 (id 5282)
$ws.Range("H46").Value = 1666
$ws.Range("I46").Value = 2000
$ws.Range("J46").Value = 998
$ws.Range("K46").Value = 2000
$ws.Range("L46").Value = 998
$ws.Range("M46").Value = -1812
$ws.Range("N46").Value = -1374

# Row 68 (id 12563)
$ws.Range("H68").Value = 2374.3333
$ws.Range("I68").Value = 1898.3334
$ws.Range("J68").Value = 2850.3333
$ws.Range("K68").Value = 1898.3334
$ws.Range("L68").Value = 2850.3333
$ws.Range("M68").Value = -1149.3334
$ws.Range("N68").Value = -4348.3333

# Row 71 (id 12563)
$ws.Range("H71").Value = 2374.3333
$ws.Range("I71").Value = 1898.3334
$ws.Range("J71").Value = 2850.3333
$ws.Range("K71").Value = 9491.666999999999
$ws.Range("L71").Value = 14251.6665
$ws.Range("M71").Value = -5747.666999999999
$ws.Range("N71").Value = -21739.6665

# Row 126 (id 36249)
$ws.Range("H126").Value = 2219.6538
$ws.Range("I126").Value = 2054.1333
$ws.Range("J126").Value = 2445.3635
$ws.Range("K126").Value = 6162.3999
$ws.Range("L126").Value = 7336.0905
$ws.Range("M126").Value = -3692.3999
$ws.Range("N126").Value = -12276.0905

# Row 132 (id 44058)
$ws.Range("H132").Value = 15698.75
$ws.Range("I132").Value = 4733.1665
$ws.Range("J132").Value = 22278.1
$ws.Range("K132").Value = 14199.4995
$ws.Range("L132").Value = 66834.29999999999
$ws.Range("M132").Value = -11669.4995
$ws.Range("N132").Value = -71894.29999999999

$ws = $wb.Worksheets.Item("WVR")
# Row 107 (id 27746)
$ws.Range("H107").Value = 491.3889
$ws.Range("I107").Value = 428.75
$ws.Range("J107").Value = 616.6667
$ws.Range("K107").Value = 1286.25
$ws.Range("L107").Value = 1850.0001
$ws.Range("M107").Value = 633.75
$ws.Range("N107").Value = -5690.0001

# Row 113 (id 27752)
$ws.Range("H113").Value = 260.76923
$ws.Range("I113").Value = 294
$ws.Range("J113").Value = 240
$ws.Range("K113").Value = 882
$ws.Range("L113").Value = 720
$ws.Range("M113").Value = 1288
$ws.Range("N113").Value = -5060

# Row 122 (id 36208)
$ws.Range("H122").Value = 43215.25
$ws.Range("I122").Value = 72552.71000000001
$ws.Range("J122").Value = 2142.8
$ws.Range("K122").Value = 217658.13
$ws.Range("L122").Value = 6428.400000000001
$ws.Range("M122").Value = -215208.13
$ws.Range("N122").Value = -11328.4

# Row 136 (id 44031)
$ws.Range("H136").Value = 29631914
$ws.Range("I136").Value = 34484492
$ws.Range("J136").Value = 20836614
$ws.Range("K136").Value = 103453476
$ws.Range("L136").Value = 62509842
$ws.Range("M136").Value = -103450926
$ws.Range("N136").Value = -62514942
